$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old column C (Tasa con 25kV) entirely
$ws.Range("C1:C5").ClearContents()

# Update header row 1
$ws.Range("B1").Value = "Tasa con 35kV"

# Update header row 2
$ws.Range("B2").Value = "R(35kV)/Imp/s"

# Clear out the previous data rows (row 4 and 5) before rewriting the full data block
$ws.Range("A4:C9").ClearContents()

# Write new data rows 4-9
$data = @(
    @(0, 13080),
    @(10, 13031),
    @(20, 12485),
    @(30, 13283),
    @(40, 12750),
    @(50, 11838)
)

$r = 4
foreach ($pair in $data) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# Match the final selection to the new used range
$ws.Range("A1:B9").Select()
